# Corrige las cifras del informe PDF.
# 1) Rellena la columna "Ventas" (K) con las cifras correctas que faltaban.
# 2) Elimina las columnas M y N ("Unnamed: 12" / "Unnamed: 13", vacias/erroneas)
#    desplazando el resto de columnas de calculo (Demanda diaria, Stock minimo,
#    Stock seguridad, Stock maximo, Mes) dos posiciones a la izquierda.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Nuevos valores de "Ventas" (columna K) ---
$ventas = @{
    2  = 900000
    3  = 300000
    4  = 100000
    5  = 50000
    6  = 890000
    7  = 1000000
    8  = 1200000
    9  = 950000
    10 = 350000
    11 = 350000
    12 = 350000
    13 = 350000
    14 = 350000
}

foreach ($row in $ventas.Keys) {
    $ws.Cells.Item($row, 11).Value = $ventas[$row]
}

# --- 2) Elimina las columnas M:N (desplaza O:S -> M:Q) ---
$ws.Range("M1:N1").EntireColumn.Delete()

Write-Output "done"
